# Add 2022-Q3 data:
#  1. Insert a new row at the top of the "总计" (summary) sheet's data with the
#     2022-Q3 totals, shifting the existing quarters down by one row.
#  2. Insert a brand-new worksheet named "2022-Q3" (right after "总计") holding
#     the per-fund breakdown for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" sheet: shift existing rows down and write the new 2022-Q3 row.
#    (Values are written as literals rather than via Rows.Insert()/read-back
#    so every cell keeps the exact type/format it should have.)
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)

$totalSheet.Cells.Item(7, 1).Value = 5
$totalSheet.Cells.Item(7, 2).Value = "2021-Q1"
$totalSheet.Cells.Item(7, 3).Value = 5
$totalSheet.Cells.Item(7, 4).Value = 1.01

$totalSheet.Cells.Item(6, 1).Value = 4
$totalSheet.Cells.Item(6, 2).Value = "2021-Q2"
$totalSheet.Cells.Item(6, 3).Value = 1
$totalSheet.Cells.Item(6, 4).Value = 0.05

$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(5, 2).Value = "2021-Q3"
$totalSheet.Cells.Item(5, 3).Value = 4
$totalSheet.Cells.Item(5, 4).Value = 0.87

$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(4, 2).Value = "2021-Q4"
$totalSheet.Cells.Item(4, 3).Value = 16
$totalSheet.Cells.Item(4, 4).Value = 4.66

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(3, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(3, 3).Value = 8
$totalSheet.Cells.Item(3, 4).Value = 2.22

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 12
$totalSheet.Cells.Item(2, 4).Value = 3.77

# New row-2 and row-7 "index column" cells must carry the same bold/bordered
# style that every other A-column cell in this sheet already uses (row 7 is a
# brand-new row that never existed before, so it has no style at all yet).
$totalSheet.Cells.Item(3, 1).Copy()
$totalSheet.Cells.Item(2, 1).PasteSpecial(-4122)  # xlPasteFormats
$totalSheet.Cells.Item(7, 1).PasteSpecial(-4122)  # xlPasteFormats
$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(7, 1).Value = 5

# ---------------------------------------------------------------------------
# 2. Brand-new "2022-Q3" worksheet, inserted right after "总计".
# ---------------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($beforeSheet)
$newSheet.Name = "2022-Q3"

# Re-fetch by name: the COM layer here re-targets any existing sheet handle
# (like $beforeSheet) onto the freshly-inserted sheet once Add() runs, so we
# look sheets up by name afterwards instead of trusting old references.
$q1Sheet = $wb.Worksheets.Item("2022-Q1")
$q3Sheet = $wb.Worksheets.Item("2022-Q3")

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3Sheet.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Fund rows: code, name, size, total stock position, position ratio,
# held-value(亿元), position rank. All but the name and the rank are text in
# the source data, including numeric-looking strings like "9.00" - the
# leading "'" forces text storage, and resetting .Style to "Normal" afterwards
# clears the stray "quote-prefix" formatting flag that the apostrophe leaves
# behind, while keeping the value itself stored as text.
$fundRows = @(
    @("519702", "交银趋势优先混合A",         "83.94", "81.61", "2.18", "1.8299", 10),
    @("014038", "交银启诚混合A",             "24.82", "81.41", "2.72", "0.6751", 8),
    @("001128", "宝盈新兴产业灵活配置混合A",  "9.00",  "91.66", "3.19", "0.2871", 9),
    @("001487", "宝盈优势产业灵活配置混合A",  "10.11", "91.85", "2.82", "0.2851", 9),
    @("013430", "交银趋势优先混合C",         "12.61", "81.61", "2.18", "0.2749", 10),
    @("014039", "交银启诚混合C",             "7.63",  "81.41", "2.72", "0.2075", 8),
    @("012771", "宝盈优势产业灵活配置混合C",  "3.62",  "91.85", "2.82", "0.1021", 9),
    @("011404", "融通鑫新成长混合C",         "1.75",  "94.07", "2.93", "0.0513", 8),
    @("012815", "宝盈新兴产业灵活配置混合C",  "1.31",  "91.66", "3.19", "0.0418", 9),
    @("011403", "融通鑫新成长混合A",         "0.39",  "94.07", "2.93", "0.0114", 8),
    @("003855", "汇安丰华灵活配置混合C",     "0.19",  "45.55", "1.84", "0.0035", 10),
    @("003854", "汇安丰华灵活配置混合A",     "0.00",  "45.55", "1.84", $null,    10)
)

for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $row = $i + 2
    $data = $fundRows[$i]

    $q3Sheet.Cells.Item($row, 1).Value = $i          # A: running index (number)
    $q3Sheet.Cells.Item($row, 2).Value = "'" + $data[0]   # B: fund code (text)
    $q3Sheet.Cells.Item($row, 2).Style = "Normal"
    $q3Sheet.Cells.Item($row, 3).Value = $data[1]         # C: fund name (text, safe as-is)
    $q3Sheet.Cells.Item($row, 4).Value = "'" + $data[2]   # D: fund size (text)
    $q3Sheet.Cells.Item($row, 4).Style = "Normal"
    $q3Sheet.Cells.Item($row, 5).Value = "'" + $data[3]   # E: total stock position (text)
    $q3Sheet.Cells.Item($row, 5).Style = "Normal"
    $q3Sheet.Cells.Item($row, 6).Value = "'" + $data[4]   # F: position ratio (text)
    $q3Sheet.Cells.Item($row, 6).Style = "Normal"
    if ($null -ne $data[5]) {
        $q3Sheet.Cells.Item($row, 7).Value = "'" + $data[5]   # G: held value (text)
        $q3Sheet.Cells.Item($row, 7).Style = "Normal"
    }
    $q3Sheet.Cells.Item($row, 8).Value = $data[6]         # H: position rank (number)
}
# The very last fund's "held value" is a genuine 0 number in the source data
# (not "0.0000" text) - set separately since it's the one row left out above.
$q3Sheet.Cells.Item(13, 7).Value = 0

# Header row + index column get the same bold/bordered style used on every
# other quarterly sheet.
$q1Sheet.Range("B1:H1").Copy()
$q3Sheet.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats
$q1Sheet.Range("A2").Copy()
$q3Sheet.Range("A2:A13").PasteSpecial(-4122)  # xlPasteFormats
